$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: "GoogleN" -> "IntelN" for rows 2-16
for ($i = 1; $i -le 15; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = "Intel$i"
}

# Column B, row 16 (last row) gets its new value first: "Productbased15" -> "Servicebased15"
$ws.Cells.Item(16, 2).Value = "Servicebased15"

# Column B: "ProductbasedN" -> "ServicebasedN" for rows 2-15
for ($i = 1; $i -le 14; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = "Servicebased$i"
}

# Update the selection to C2 (single cell) as shown in the diff
$ws.Range("C2").Select()
